$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.92
$wsSummary.Range("B4").Value = -0.08
$wsSummary.Range("B5").Value = -0.18
$wsSummary.Range("B6").Value = 9
$wsSummary.Range("B7").Value = 3
$wsSummary.Range("B9").Value = 33.33

# --- Sheet: Strategy Status ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 99.92
$wsStrategy.Range("D4").Value = 9
$wsStrategy.Range("E4").Value = -0.08
$wsStrategy.Range("F4").Value = -0.08
$wsStrategy.Range("G4").Value = 33.33

# --- Sheet: All Trades ---
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G10").Value = 0.95
$wsTrades.Range("H10").Value = "CLOSED"
$wsTrades.Range("I10").Value = 2.1505
$wsTrades.Range("J10").Value = 0.02
$wsTrades.Range("K10").Value = 99.92
$wsTrades.Range("P10").Value = "early_exit"
$wsTrades.Range("Q10").Value = 0.13

# --- Sheet: MarketMaking ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G10").Value = 0.95
$wsMM.Range("H10").Value = "CLOSED"
$wsMM.Range("I10").Value = 2.1505
$wsMM.Range("J10").Value = 0.02
$wsMM.Range("K10").Value = 99.92
$wsMM.Range("P10").Value = "early_exit"
$wsMM.Range("Q10").Value = 0.13
